# Update for first draft
# Rewrites the cd_rcv.xlsx-style results table: split each of the three
# (State Based / Non State / One Sided) columns into a 'mean' + 'std' pair,
# rename CART -> DTREE, drop NB (and its data row), and refresh all values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): relabel existing columns, add 3 new std columns ---
$ws.Range("B1").Value2 = "Algorithm"
$ws.Range("C1").Value2 = "State Based mean"
$ws.Range("D1").Value2 = "State Based std"
$ws.Range("E1").Value2 = "Non State mean"
$ws.Range("F1").Value2 = "Non State std"
$ws.Range("G1").Value2 = "One Sided mean"
$ws.Range("H1").Value2 = "One Sided std"

# Copy the existing header style (bold, centered, bordered) onto the new
# F1:H1 header cells so they match B1:E1 exactly.
$ws.Range("B1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows 2-8: refresh algorithm labels + mean values, add std columns ---
$ws.Range("A2").Value2 = 0
$ws.Range("B2").Value2 = "LR"
$ws.Range("C2").Value2 = 0.8971001896333755
$ws.Range("D2").Value2 = 0.02683384970588557
$ws.Range("E2").Value2 = 0.8791858678955453
$ws.Range("F2").Value2 = 0.03489014928222599
$ws.Range("G2").Value2 = 0.8933277731442869
$ws.Range("H2").Value2 = 0.0196568647042139

$ws.Range("A3").Value2 = 1
$ws.Range("B3").Value2 = "LDA"
$ws.Range("C3").Value2 = 0.9086678255372946
$ws.Range("D3").Value2 = 0.02927036970641505
$ws.Range("E3").Value2 = 0.9065796210957501
$ws.Range("F3").Value2 = 0.03875422606649659
$ws.Range("G3").Value2 = 0.914303586321935
$ws.Range("H3").Value2 = 0.01430339179504787

$ws.Range("A4").Value2 = 2
$ws.Range("B4").Value2 = "KNN"
$ws.Range("C4").Value2 = 0.9326011378002528
$ws.Range("D4").Value2 = 0.01581265494325974
$ws.Range("E4").Value2 = 0.9323092677931388
$ws.Range("F4").Value2 = 0.02962679306039942
$ws.Range("G4").Value2 = 0.9334278565471225
$ws.Range("H4").Value2 = 0.01430079217560514

$ws.Range("A5").Value2 = 3
$ws.Range("B5").Value2 = "DTREE"
$ws.Range("C5").Value2 = 0.769326801517067
$ws.Range("D5").Value2 = 0.03675448364489299
$ws.Range("E5").Value2 = 0.6748079877112134
$ws.Range("F5").Value2 = 0.0496699235415331
$ws.Range("G5").Value2 = 0.736555462885738
$ws.Range("H5").Value2 = 0.03554031660914911

$ws.Range("A6").Value2 = 4
$ws.Range("B6").Value2 = "RTREE"
$ws.Range("C6").Value2 = 0.9077749683944374
$ws.Range("D6").Value2 = 0.02953923664450255
$ws.Range("E6").Value2 = 0.9001280081925243
$ws.Range("F6").Value2 = 0.03877553038484891
$ws.Range("G6").Value2 = 0.9015012510425354
$ws.Range("H6").Value2 = 0.02681678081885578

$ws.Range("A7").Value2 = 5
$ws.Range("B7").Value2 = "XTREE"
$ws.Range("C7").Value2 = 0.8500790139064476
$ws.Range("D7").Value2 = 0.03751796611568015
$ws.Range("E7").Value2 = 0.8549923195084486
$ws.Range("F7").Value2 = 0.03768108957928044
$ws.Range("G7").Value2 = 0.8413678065054212
$ws.Range("H7").Value2 = 0.02057022353591244

$ws.Range("A8").Value2 = 6
$ws.Range("B8").Value2 = "SVM"
$ws.Range("C8").Value2 = 0.8917746523388116
$ws.Range("D8").Value2 = 0.02872272441916114
$ws.Range("E8").Value2 = 0.8727342549923195
$ws.Range("F8").Value2 = 0.04819047903038755
$ws.Range("G8").Value2 = 0.89790658882402
$ws.Range("H8").Value2 = 0.02025615884426553

# --- Row 9 (old 'NB' row) no longer exists in the new layout; remove it ---
$ws.Rows("9").Delete()

Write-Host "Applied cd_rcv update: mean/std split, DTREE rename, NB row removed."
